# Generate Report for Handoff
# Updates the localization-status workbook so that the "b.md" row reflects
# that it is now ready for handoff (instead of "Handed back: in sync with
# en-US"), with a fresh handoff file / timestamp and an error detail message
# explaining that the handback file version is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68dd6ca6c64fa41d8f8bf16b75b39151df6abb04/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41380c5b7356355012fe290d666745ab2ef37a7c/e2e/b.md."

# --- Overview sheet: b.md row (row 3) -------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 22:38:14"

# --- zh-cn sheet: b.md row (row 3) -----------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 22:38:10"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: b.md row (row 3) -----------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 22:38:14"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
